$d = $word.ActiveDocument

# Avoid Word's "AutoFormat As You Type" turning straight quotes into curly
# ("smart") quotes when we type literal double-quote characters later on.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function Find-ParagraphIndex($substr) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Title: "Alexandru " + "Ardelean" (two runs) -> "Alexandru Ardelean" (one)
# ---------------------------------------------------------------------------
# A plain Range.Text assignment is a no-op here because the visible text
# does not actually change, so the two runs would be left untouched. Using
# Find & Replace instead forces Word to rebuild the run, merging both
# original runs into a single one.
$found = $d.Content.Find.Execute("Alexandru Ardelean", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Alexandru Ardelean", 2)

# ---------------------------------------------------------------------------
# 2) "Linux kernel contributions..." paragraph split into 3 paragraphs
# ---------------------------------------------------------------------------
$linuxIdx = Find-ParagraphIndex("Linux kernel contributions")
$linuxPara = $d.Paragraphs($linuxIdx).Range
$oldText = "Linux kernel contributions: roughly 1400+ (code contributions and reviews)."
$prefix = "Linux kernel contributions: roughly "
$tailStart = $linuxPara.Start + $prefix.Length
$tailEnd = $linuxPara.End - 1
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Text = "~430 authored"

# Insert a brand-new paragraph right after it, inheriting the same
# paragraph/run formatting, then fill it with the git command line.
$linuxPara2 = $d.Paragraphs($linuxIdx).Range
$newPara = $linuxPara2.InsertParagraphAfter()
$gitParaRange = $d.Paragraphs($linuxIdx + 1).Range
$gitInsertPoint = $d.Range($gitParaRange.Start, $gitParaRange.Start)
$gitInsertPoint.Text = '  (  git log --oneline --author="Alexandru Ardelean"  | wc -l )'

# ---------------------------------------------------------------------------
# 3) "More recently, ..." paragraph: drop the power-consumption aside
# ---------------------------------------------------------------------------
$moreIdx = Find-ParagraphIndex("More recently, I started doing Computer Vision")
$found = $d.Content.Find.Execute(" (taking up 2-3 Watts of power)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# 4) Skill-set table: column widths shift by one twip (0.05pt) each
# ---------------------------------------------------------------------------
$skillsTable = $d.Tables(1)
$skillsTable.Columns(1).Width = 112.4
$skillsTable.Columns(2).Width = 153.1
